$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new key/value (building-block) pairs right after the existing
# data (rows 1-315 already populated). Writing order below reproduces the
# author's actual entry sequence (descriptive text for the first new row,
# then the batch of new z-codes, then the batch of remaining descriptions,
# then the final z-code) so the shared-string table grows in the same
# order as the recorded edit.
$ws.Range("B316").Value = "정사영된 삼각형의 각 꼭짓점을 파악해서 정사영된 삼각형 넓이의 최댓값을 구합니다."

$ws.Range("A316").Value = "z3006"
$ws.Range("A317").Value = "z3007"
$ws.Range("A318").Value = "z3008"
$ws.Range("A319").Value = "z3009"

$ws.Range("B317").Value = "피타고라스 정리를 이용해서 정사영의 넓이를 최대화 시키는 삼각형 `$\mathrm{PQR}`$의 세변의 길이를 구합니다."
$ws.Range("B318").Value = "이등변삼각형 `$\mathrm{PQR}`$의 넓이를 구합니다."
$ws.Range("B319").Value = "정사영 전과 후의 두 넓이의 비율을 이용해서 두 평면이 이루는 예각에 대한 코사인 값을 구합니다."
$ws.Range("B320").Value = "코사인 값을 이용해서 정사영의 넓이를 구합니다."

$ws.Range("A320").Value = "z3010"

# Match the author's final cursor/viewport position recorded in the diff
# (selection moved to C313 while scrolled down near row 295).
$ws.Range("C313").Select()
